# Update the cryptos list (GitHub Actions scheduled data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.204.63"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "2.068.08"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Formula = "'251.03"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Formula = "'0.677"
$ws.Range("E6").Value = "  +4.26%  "
$ws.Range("D7").Formula = "'62.32"
$ws.Range("E7").Value = "  +25.71%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Formula = "'61.16"
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("D10").Formula = "'0.383"
$ws.Range("E10").Value = "  +4.34%  "
$ws.Range("D11").Formula = "'0.0805"
$ws.Range("E11").Value = "  +9.77%  "
$ws.Range("D12").Formula = "'0.108"
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("D13").Formula = "'15.61"
$ws.Range("E13").Value = "  +4.55%  "
$ws.Range("D14").Value = "2.368.00"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Formula = "'0.825"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Formula = "'5.41"
$ws.Range("E16").Value = "  +7.52%  "
$ws.Range("D17").Value = "2.064.89"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "37.174.24"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").Formula = "'75.21"
$ws.Range("E19").Value = "  +5.00%  "
$ws.Range("E20").Value = "  +14.23%  "
$ws.Range("D21").Formula = "'14.78"
$ws.Range("E21").Value = "  +12.97%  "
$ws.Range("D22").Formula = "'5.43"
$ws.Range("E22").Value = "  +5.79%  "
$ws.Range("D23").Formula = "'240.00"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Formula = "'2.43"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Formula = "'171.87"
$ws.Range("E26").Value = "  +2.56%  "
$ws.Range("D27").Formula = "'9.25"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").Formula = "'20.39"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").Formula = "'2.02"
$ws.Range("E29").Value = "  +1.52%  "
$ws.Range("D30").Formula = "'0.126"
$ws.Range("E30").Value = "  +3.71%  "
$ws.Range("D33").Formula = "'0.0638"
$ws.Range("E33").Value = "  +6.54%  "
$ws.Range("D34").Formula = "'4.44"
$ws.Range("E34").Value = "  +10.27%  "
$ws.Range("D35").Formula = "'0.0889"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("D36").Formula = "'1.00"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").Formula = "'2.31"
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("E38").Value = "  -3.49%  "
$ws.Range("D39").Formula = "'0.111"
$ws.Range("E39").Value = "  +28.58%  "
$ws.Range("E40").Value = "  +3.92%  "
$ws.Range("D41").Formula = "'18.76"
$ws.Range("E41").Value = "  +7.57%  "
$ws.Range("D42").Formula = "'0.0227"
$ws.Range("E42").Value = "  +2.95%  "
$ws.Range("D43").Formula = "'1.16"
$ws.Range("E43").Value = "  +1.78%  "
$ws.Range("D44").Formula = "'98.28"
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("D45").Formula = "'4.28"
$ws.Range("E45").Value = "  +32.87%  "
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("D49").Value = "1.307.92"
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("E51").Value = "  +2.07%  "

# Row 31 becomes ImmutableX (was Filecoin), Row 32 becomes Filecoin (was ImmutableX)
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Formula = "'1.09"
$ws.Range("E31").Value = "  +1.44%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Formula = "'4.65"
$ws.Range("E32").Value = "  +4.84%  "

# Row 47 becomes RenderToken (was THORChain), Row 48 becomes THORChain (was RenderToken)
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Formula = "'2.54"
$ws.Range("E47").Value = "  +13.85%  "

$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Formula = "'4.51"
$ws.Range("E48").Value = "  +16.21%  "
